# Apply the "Several rows supplemented" update to the ChartSettings sheet.
$wb = $excel.ActiveWorkbook
$wsSettings = $wb.Worksheets.Item("ChartSettings")

# --- Update existing row 2 values on ChartSettings ---
$wsSettings.Cells.Item(2,2).Value = 45541      # B2 Date_Change -> 2024-09-06
$wsSettings.Cells.Item(2,19).Value = 0.12      # S2 Grid_Bottom

# --- New HeatNeed rows to append (row, id, dateSerial, S-value) ---
$newRows = @(
    @{ Row = 3; Id = "HeatNeed.02"; Date = 45436; S = 0.2  },
    @{ Row = 4; Id = "HeatNeed.03"; Date = 45541; S = 0.25 },
    @{ Row = 5; Id = "HeatNeed.04"; Date = 45541; S = 0.3  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $wsSettings.Cells.Item($row,1).Value = $r.Id          # A - ID_Settings
    $wsSettings.Cells.Item($row,2).Value = $r.Date         # B - Date_Change
    $wsSettings.Cells.Item($row,2).NumberFormat = "m/d/yy" # keep the date number format (style 2)

    $wsSettings.Cells.Item($row,4).Value = "Energy Balance Building"       # D - ChartTitle_ENG
    $wsSettings.Cells.Item($row,5).Value = "Energiebilanz Gebäude"         # E - ChartTitle_GER
    $wsSettings.Cells.Item($row,7).Value = "Heat transfer during heating season"   # G - ChartSubTitle_ENG
    $wsSettings.Cells.Item($row,8).Value = "Wärmestrom während der Heizperiode"    # H - ChartSubTitle_GER
    $wsSettings.Cells.Item($row,13).Value = "Heat transfer related to reference area [kWh/(m²a)]"        # M - AxisTitle_y_ENG
    $wsSettings.Cells.Item($row,14).Value = "Wärmestrom bezogen auf Referenzfläche [kWh/(m²a)]"           # N - AxisTitle_y_GER

    $wsSettings.Cells.Item($row,16).Value = 20    # P - FontSize
    $wsSettings.Cells.Item($row,17).Value = 20    # Q - FontSize_Legend

    $wsSettings.Cells.Item($row,18).Value = 0.12   # R - Grid_Top
    $wsSettings.Cells.Item($row,18).NumberFormat = "0%"
    $wsSettings.Cells.Item($row,19).Value = $r.S   # S - Grid_Bottom
    $wsSettings.Cells.Item($row,19).NumberFormat = "0%"
    $wsSettings.Cells.Item($row,20).Value = 0.1    # T - Grid_Left
    $wsSettings.Cells.Item($row,20).NumberFormat = "0%"
    $wsSettings.Cells.Item($row,21).Value = 0.1    # U - Grid_Right
    $wsSettings.Cells.Item($row,21).NumberFormat = "0%"

    $wsSettings.Cells.Item($row,22).Value = 0      # V - AxisMin_y
    $wsSettings.Cells.Item($row,23).Value = 300    # W - AxisMax_y
    $wsSettings.Cells.Item($row,25).Value = 5      # Y - AxisMinInterval_y
    $wsSettings.Cells.Item($row,26).Value = 50     # Z - AxisMaxInterval_y

    $wsSettings.Rows.Item($row).RowHeight = 28.8
}

# --- View state: ChartSettings becomes the active/selected sheet ---
$wsSettings.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 1
$wsSettings.Range("S5").Select()
